# REVER_DailyTracker_MONISHA.xlsx - Add files via upload
# Updates the FEB-2021 sheet (6th sheet / tab) with new daily task entries
# for 18-Feb to 24-Feb-2021, moves the "Week Off" block from D24:D25 to
# D21:D22, and fills in D24/D25 with actual task info.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEB-2021")

# ---------------------------------------------------------------------
# Row 19 (18-Feb-2021): Hayyai - B2B / New Layout ... Modification / 100% / Completed
# ---------------------------------------------------------------------
$ws.Range("C17:F17").Copy()
$ws.Range("C19:F19").PasteSpecial(-4122)
$ws.Range("C19").Value = "Hayyai - B2B"
$ws.Range("D19").Value = "New Layout -Service Job Management - Acknowledgement & Engineer Assign, create Job - 2 screens Modification"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = "Completed"
$ws.Rows.Item(19).RowHeight = 43.2

# ---------------------------------------------------------------------
# Row 20 (19-Feb-2021): Hayyai - B2B / Modification in Registration layout / 100% / Completed
# ---------------------------------------------------------------------
$ws.Range("C17:F17").Copy()
$ws.Range("C20:F20").PasteSpecial(-4122)
$ws.Range("C20").Value = "Hayyai - B2B"
$ws.Range("D20").Value = "Modification in Registration layout"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = "Completed"

# ---------------------------------------------------------------------
# Rows 21-22 (20/21-Feb-2021): the merged "Week Off" block moves here
# (it used to live at D24:D25).
# ---------------------------------------------------------------------
$ws.Range("D24:D25").UnMerge()
$ws.Range("D21:D22").Merge()
$ws.Range("D24:D25").Copy()
$ws.Range("D21:D22").PasteSpecial(-4122)
$ws.Range("D21").Value = "Week Off"

# ---------------------------------------------------------------------
# Row 23 (22-Feb-2021): Hayyai - B2B / Service Job Management - Search Job screen new layout / 80% / WIP
# ---------------------------------------------------------------------
$ws.Range("C18:F18").Copy()
$ws.Range("C23:F23").PasteSpecial(-4122)
$ws.Range("C23").Value = "Hayyai - B2B"
$ws.Range("D23").Value = "Service Job Management - Serach Job screen new layout"
$ws.Range("E23").Value = 0.8
$ws.Range("F23").Value = "WIP"
$ws.Rows.Item(23).RowHeight = 28.8

# ---------------------------------------------------------------------
# Row 24 (23-Feb-2021): Hayyai - B2B / Dashboard - Pending calls screen / 100% / Completed
# ---------------------------------------------------------------------
$ws.Range("C17:F17").Copy()
$ws.Range("C24:F24").PasteSpecial(-4122)
$ws.Range("C24").Value = "Hayyai - B2B"
$ws.Range("D24").Value = "Dashboard - Pending calls screen"
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = "Completed"

# ---------------------------------------------------------------------
# Row 25 (24-Feb-2021): D25 is no longer part of a merge - restore the
# plain (non-merged) wrap-text style used elsewhere in the column.
# ---------------------------------------------------------------------
$ws.Range("D26").Copy()
$ws.Range("D25").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# View state: scroll position / active cell on the FEB-2021 sheet, and
# the workbook window geometry.
# ---------------------------------------------------------------------
$ws.Range("D26").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

$excel.ActiveWindow.Left = 4080
$excel.ActiveWindow.Top = 756
$excel.ActiveWindow.Width = 17280
$excel.ActiveWindow.Height = 11604
